# Apply Coinranking crypto price/volume refresh (GitHub Actions update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.126.50"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").Value = "3.101.75"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'573.62"
$ws.Range("E5").Value = "  -0.97%  "
$ws.Range("D6").Value = "'178.03"
$ws.Range("E6").Value = "  +3.77%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "3.099.03"
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("E9").Value = "  -1.05%  "
$ws.Range("E10").Value = "  -1.09%  "
$ws.Range("D11").Value = "'0.151"
$ws.Range("E11").Value = "  -0.20%  "
$ws.Range("D12").Value = "'0.467"
$ws.Range("E12").Value = "  -1.89%  "
$ws.Range("E13").Value = "  -2.01%  "
$ws.Range("D14").Value = "'36.16"
$ws.Range("E14").Value = "  -0.98%  "
$ws.Range("E15").Value = "  +0.33%  "
$ws.Range("D16").Value = "3.622.04"
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("D17").Value = "67.081.80"
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("E18").Value = "  -0.74%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "'16.74"
$ws.Range("E19").Value = "  +0.50%  "
$ws.Range("B20").Value = "WrappedEther"
$ws.Range("C20").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D20").Value = "3.105.10"
$ws.Range("E20").Value = "  +0.53%  "
$ws.Range("D21").Value = "'490.00"
$ws.Range("E21").Value = "  +0.59%  "
$ws.Range("D22").Value = "'7.73"
$ws.Range("E22").Value = "  -0.72%  "
$ws.Range("D23").Value = "'0.686"
$ws.Range("E23").Value = "  -1.39%  "
$ws.Range("D24").Value = "'83.44"
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("E25").Value = "  +0.44%  "
$ws.Range("E26").Value = "  -2.97%  "
$ws.Range("E27").Value = "  -2.24%  "
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("E29").Value = "  +1.35%  "
$ws.Range("D30").Value = "'2.30"
$ws.Range("E30").Value = "  -0.89%  "
$ws.Range("D31").Value = "'2.59"
$ws.Range("E31").Value = "  -2.21%  "
$ws.Range("D32").Value = "'28.08"
$ws.Range("E32").Value = "  -0.55%  "
$ws.Range("E33").Value = "  -0.89%  "
$ws.Range("D34").Value = "0.0₃0942"
$ws.Range("E34").Value = "  +0.19%  "
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("B36").Value = "Arweave"
$ws.Range("C36").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D36").Value = "'47.37"
$ws.Range("E36").Value = "  +2.01%  "
$ws.Range("B37").Value = "Mantle"
$ws.Range("C37").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D37").Value = "'0.947"
$ws.Range("E37").Value = "  -1.92%  "
$ws.Range("D38").Value = "'5.57"
$ws.Range("E38").Value = "  -3.26%  "
$ws.Range("E39").Value = "  +2.52%  "
$ws.Range("E40").Value = "  +0.46%  "
$ws.Range("D41").Value = "'49.14"
$ws.Range("E41").Value = "  -1.20%  "
$ws.Range("D42").Value = "'0.123"
$ws.Range("E42").Value = "  +0.14%  "
$ws.Range("D43").Value = "'8.24"
$ws.Range("E43").Value = "  -2.14%  "
$ws.Range("E44").Value = "  +5.37%  "
$ws.Range("D45").Value = "2.799.30"
$ws.Range("E45").Value = "  +0.38%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "'0.0345"
$ws.Range("E46").Value = "  -1.03%  "
$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").Value = "'369.14"
$ws.Range("E47").Value = "  -2.93%  "
$ws.Range("D48").Value = "'135.58"
$ws.Range("E48").Value = "  +0.38%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("E50").Value = "  +3.09%  "
$ws.Range("D51").Value = "'2.28"
$ws.Range("E51").Value = "  +4.66%  "

Write-Output "Applied cryptos update"
